$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update A12: model id
$ws.Range("A12").Value = "BIOMD0000000991"

# Update C12: new odes text block
$newOdes = "odes = [`n" +
"    sympy.Eq(S(t).diff(t), - beta_c * (alpha * A(t) + I(t)) / (Nh - ID_param) * S(t)),`n" +
"    sympy.Eq(E(t).diff(t), beta_c * (alpha * A(t) + I(t)) / (Nh - ID_param) * S(t) - sigma * E(t)),`n" +
"    sympy.Eq(A(t).diff(t), nu * sigma * E(t) - (theta + gamma_a) * A(t)),`n" +
"    sympy.Eq(I(t).diff(t), (1 - nu) * sigma * E(t) - (psi + gamma_O + dO) * I(t)),`n" +
"    sympy.Eq(ID(t).diff(t), theta * A(t) + psi * I(t) - (gamma_i + dD) * ID(t)),`n" +
"    sympy.Eq(R(t).diff(t), gamma_i * ID(t) + gamma_a * A(t) + gamma_O * I(t))`n" +
"]"

$ws.Range("C12").Value = $newOdes
